# Updated symbol list on Tue Dec 27 20:50:48 UTC 2022 with GitHub Actions
# Refresh the cryptocurrency price / volume snapshot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    # Prefix with an apostrophe so a numeric-looking string ("245.67") is
    # stored as text rather than coerced into a Number cell, matching the
    # original inline-string cell type. Restore the default style afterwards
    # so no stray numeric/text style index gets attached to the cell.
    $rng.Value = "'" + $Value
    $rng.Style = "Normal"
}

# Price column (D) updates
Set-TextValue "D2" "245.67"
Set-TextValue "D3" "24.00"
Set-TextValue "D4" "5.369"
Set-TextValue "D5" "0.05826"
Set-TextValue "D6" "6.460"
Set-TextValue "D7" "3.360"
Set-TextValue "D8" "0.8098"
Set-TextValue "D9" "0.9198"
Set-TextValue "D10" "0.1409"
Set-TextValue "D11" "0.07356"
Set-TextValue "D12" "0.03112"
Set-TextValue "D13" "0.03051"
Set-TextValue "D14" "0.09373"
Set-TextValue "D15" "3.854"
Set-TextValue "D16" "0.001561"
Set-TextValue "D17" "0.04695"
Set-TextValue "D18" "0.0005991"
Set-TextValue "D19" "0.006176"
Set-TextValue "D20" "0.001246"
Set-TextValue "D21" "0.004693"
Set-TextValue "D22" "0.00008801"
Set-TextValue "D23" "3.596"
Set-TextValue "D25" "0.3184"
Set-TextValue "D28" "0.0002350"
Set-TextValue "D41" "0.006441"
Set-TextValue "D42" "0.1065"
Set-TextValue "D43" "0.003200"
Set-TextValue "D44" "0.008527"
Set-TextValue "D45" "0.00005254"
Set-TextValue "D47" "0.6861"

# Volume(1h) column (E) label updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
